# Reproduces the commit "Updating test files to match the current format in beta"
# for the optimization_parameters / network_weights sheets (and the resulting
# workbook-level active-tab bookkeeping).

$wb = $excel.ActiveWorkbook

# --- optimization_parameters sheet -----------------------------------------
$opt = $wb.Worksheets.Item("optimization_parameters")

# Row 1 used to repeat the "value" header across C1:F1 - only A1:B1 remain.
$opt.Range("C1:F1").ClearContents()

# The "Model" parameter row is renamed to "production_function" ...
$opt.Range("A8").Value = "production_function"

# ... and gets a new sibling row right below it for the L-curve toggle.
$opt.Rows.Item(9).Insert()
$opt.Range("A9").Value = "L_curve"
$opt.Range("B9").Value = 0

# The old "Deletion" row (originally row 16, now shifted to row 17 because of
# the insert above) is dropped entirely.
$opt.Rows.Item(17).Delete()

# --- network_weights sheet ---------------------------------------------------
# This used to be the selected tab; selection/content itself is unchanged,
# only the "currently active" flag moves to optimization_parameters below.
$netw = $wb.Worksheets.Item("network_weights")
$netw.Activate()
$netw.Range("B2:E5").Select()

# The optimization_parameters sheet becomes the active tab, with C1:F1
# highlighted (the cells that were just cleared) - activated last so it ends
# up as the workbook's active sheet/tab.
$opt.Activate()
$opt.Range("C1:F1").Select()
